$wb = $excel.ActiveWorkbook

# --- Capture the header style (bold, centered, bordered) from the existing
# 'Bus' sheet header row (A1) before we touch that sheet's data, so that all
# new header rows across the new sheets share the exact same style. ---
$origBusWs = $wb.Worksheets.Item("Bus")
$styleSrc = $origBusWs.Range("A1")

# --- Add the four new sheets (Load, Shunt, Line, Transformer) positioned
# right after 'Voltage Source' and before 'Bus' / 'Switch'. ---
$afterSheet = $wb.Worksheets.Item("Voltage Source")
$newSheetNames = @("Load", "Shunt", "Line", "Transformer")
foreach ($name in $newSheetNames) {
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
    $newSheet.Name = $name
    $afterSheet = $newSheet
}

# ===== Populate 'Load' sheet (dimension A11:R33) =====
$loadData = New-Object 'object[,]' 23,18
$loadData[0,0] = "Positive-Sequence Constant Imepedance Load"
$loadData[1,0] = "ID"
$loadData[1,1] = "Status"
$loadData[1,2] = "Bus"
$loadData[1,3] = "P (MW)"
$loadData[1,4] = "Q (MVAr)"
$loadData[2,0] = "End of Positive Sequence Constant Imepedance Load"
$loadData[4,0] = "Positive-Sequence Constant Power Load"
$loadData[5,0] = "ID"
$loadData[5,1] = "Status"
$loadData[5,2] = "Bus"
$loadData[5,3] = "P (MW)"
$loadData[5,4] = "Q (MVAr)"
$loadData[6,0] = "End of Positive Sequence Constant Power Load"
$loadData[8,0] = "Positive-Sequence Constant Current Load"
$loadData[9,0] = "ID"
$loadData[9,1] = "Status"
$loadData[9,2] = "Bus"
$loadData[9,3] = "P (MW)"
$loadData[9,4] = "Q (MVAr)"
$loadData[10,0] = "End of Positive Sequence Constant Current Load"
$loadData[12,0] = "Single-Phase ZIP Load"
$loadData[13,0] = "ID"
$loadData[13,1] = "Status"
$loadData[13,2] = "V (kV)"
$loadData[13,3] = "Bandwidth (pu)"
$loadData[13,4] = "Conn. type"
$loadData[13,5] = "K_z"
$loadData[13,6] = "K_i"
$loadData[13,7] = "K_p"
$loadData[13,8] = "Use initial voltage?"
$loadData[13,9] = "Bus1"
$loadData[13,10] = "P1 (kW)"
$loadData[13,11] = "Q1 (kVAr)"
$loadData[14,0] = "End of SinglePhase ZIP Load"
$loadData[16,0] = "Two-Phase ZIP Load"
$loadData[17,0] = "ID"
$loadData[17,1] = "Status"
$loadData[17,2] = "V (kV)"
$loadData[17,3] = "Bandwidth (pu)"
$loadData[17,4] = "Conn. type"
$loadData[17,5] = "K_z"
$loadData[17,6] = "K_i"
$loadData[17,7] = "K_p"
$loadData[17,8] = "Use initial voltage?"
$loadData[17,9] = "Bus1"
$loadData[17,10] = "Bus2"
$loadData[17,11] = "P1 (kW)"
$loadData[17,12] = "Q1 (kVAr)"
$loadData[17,13] = "P2 (kW)"
$loadData[17,14] = "Q2 (kVAr)"
$loadData[18,0] = "End of TwoPhase ZIP Load"
$loadData[20,0] = "Three-Phase ZIP Load"
$loadData[21,0] = "ID"
$loadData[21,1] = "Status"
$loadData[21,2] = "V (kV)"
$loadData[21,3] = "Bandwidth (pu)"
$loadData[21,4] = "Conn. type"
$loadData[21,5] = "K_z"
$loadData[21,6] = "K_i"
$loadData[21,7] = "K_p"
$loadData[21,8] = "Use initial voltage?"
$loadData[21,9] = "Bus1"
$loadData[21,10] = "Bus2"
$loadData[21,11] = "Bus3"
$loadData[21,12] = "P1 (kW)"
$loadData[21,13] = "Q1 (kVAr)"
$loadData[21,14] = "P2 (kW)"
$loadData[21,15] = "Q2 (kVAr)"
$loadData[21,16] = "P3 (kW)"
$loadData[21,17] = "Q3 (kVAr)"
$loadData[22,0] = "End of ThreePhase ZIP Load"
$loadWs = $wb.Worksheets.Item("Load")
$loadWs.Range("A11:R33").Value = $loadData
$styleSrc.Copy()
$loadWs.Range("A12:E12").PasteSpecial(-4122)
$styleSrc.Copy()
$loadWs.Range("A16:E16").PasteSpecial(-4122)
$styleSrc.Copy()
$loadWs.Range("A20:E20").PasteSpecial(-4122)
$styleSrc.Copy()
$loadWs.Range("A24:L24").PasteSpecial(-4122)
$styleSrc.Copy()
$loadWs.Range("A28:O28").PasteSpecial(-4122)
$styleSrc.Copy()
$loadWs.Range("A32:R32").PasteSpecial(-4122)

# ===== Populate 'Shunt' sheet (dimension A11:N25) =====
$shuntData = New-Object 'object[,]' 15,14
$shuntData[0,0] = "Positive Sequence Shunt"
$shuntData[1,0] = "ID"
$shuntData[1,1] = "Status"
$shuntData[1,2] = "Bus"
$shuntData[1,3] = "P (MW)"
$shuntData[1,4] = "Q (MVAr)"
$shuntData[2,0] = "End of Positive Sequence Shunt"
$shuntData[4,0] = "Single-Phase Shunt"
$shuntData[5,0] = "ID"
$shuntData[5,1] = "Status"
$shuntData[5,2] = "kV (ph-gr RMS)"
$shuntData[5,3] = "Bus1"
$shuntData[5,4] = "P1 (kW)"
$shuntData[5,5] = "Q1 (kVAr)"
$shuntData[6,0] = "End of Single-Phase Shunt"
$shuntData[8,0] = "Two-Phase Shunt"
$shuntData[9,0] = "ID"
$shuntData[9,1] = "Status1"
$shuntData[9,2] = "Status2"
$shuntData[9,3] = "kV (ph-gr RMS)"
$shuntData[9,4] = "Bus1"
$shuntData[9,5] = "Bus2"
$shuntData[9,6] = "P1 (kW)"
$shuntData[9,7] = "Q1 (kVAr)"
$shuntData[9,8] = "P2 (kW)"
$shuntData[9,9] = "Q2 (kVAr)"
$shuntData[10,0] = "End of Two-Phase Shunt"
$shuntData[12,0] = "Three-Phase Shunt"
$shuntData[13,0] = "ID"
$shuntData[13,1] = "Status1"
$shuntData[13,2] = "Status2"
$shuntData[13,3] = "Status3"
$shuntData[13,4] = "kV (ph-gr RMS)"
$shuntData[13,5] = "Bus1"
$shuntData[13,6] = "Bus2"
$shuntData[13,7] = "Bus3"
$shuntData[13,8] = "P1 (kW)"
$shuntData[13,9] = "Q1 (kVAr)"
$shuntData[13,10] = "P2 (kW)"
$shuntData[13,11] = "Q2 (kVAr)"
$shuntData[13,12] = "P3 (kW)"
$shuntData[13,13] = "Q3 (kVAr)"
$shuntData[14,0] = "End of Three-Phase Shunt"
$shuntWs = $wb.Worksheets.Item("Shunt")
$shuntWs.Range("A11:N25").Value = $shuntData
$styleSrc.Copy()
$shuntWs.Range("A12:E12").PasteSpecial(-4122)
$styleSrc.Copy()
$shuntWs.Range("A16:F16").PasteSpecial(-4122)
$styleSrc.Copy()
$shuntWs.Range("A20:J20").PasteSpecial(-4122)
$styleSrc.Copy()
$shuntWs.Range("A24:N24").PasteSpecial(-4122)

# ===== Populate 'Line' sheet (dimension A11:AA29) =====
$lineData = New-Object 'object[,]' 19,27
$lineData[0,0] = "Positive-Sequence Line"
$lineData[1,0] = "ID"
$lineData[1,1] = "Status"
$lineData[1,2] = "From bus"
$lineData[1,3] = "To bus"
$lineData[1,4] = "R (pu)"
$lineData[1,5] = "X (pu)"
$lineData[1,6] = "B (pu)"
$lineData[2,0] = "End of Positive-Sequence Line"
$lineData[4,0] = "Single-Phase Line"
$lineData[5,0] = "ID"
$lineData[5,1] = "Status"
$lineData[5,2] = "Length"
$lineData[5,3] = "From1"
$lineData[5,4] = "To1"
$lineData[5,5] = "r11 (Ohm/length_unit)"
$lineData[5,6] = "x11 (Ohm/length_unit)"
$lineData[5,7] = "b11 (uS/length_unit)"
$lineData[6,0] = "End of Single-Phase Line"
$lineData[8,0] = "Two-Phase Line"
$lineData[9,0] = "ID"
$lineData[9,1] = "Status"
$lineData[9,2] = "Length"
$lineData[9,3] = "From1"
$lineData[9,4] = "From2"
$lineData[9,5] = "To1"
$lineData[9,6] = "To2"
$lineData[9,7] = "r11 (Ohm/length_unit)"
$lineData[9,8] = "x11 (Ohm/length_unit)"
$lineData[9,9] = "r21 (Ohm/length_unit)"
$lineData[9,10] = "x21 (Ohm/length_unit)"
$lineData[9,11] = "r22 (Ohm/length_unit)"
$lineData[9,12] = "x22 (Ohm/length_unit)"
$lineData[9,13] = "b11 (uS/length_unit)"
$lineData[9,14] = "b21 (uS/length_unit)"
$lineData[9,15] = "b22 (uS/length_unit)"
$lineData[10,0] = "End of Two-Phase Line"
$lineData[12,0] = "Three-Phase Line with Full Data"
$lineData[13,0] = "ID"
$lineData[13,1] = "Status"
$lineData[13,2] = "Length"
$lineData[13,3] = "From1"
$lineData[13,4] = "From2"
$lineData[13,5] = "From3"
$lineData[13,6] = "To1"
$lineData[13,7] = "To2"
$lineData[13,8] = "To3"
$lineData[13,9] = "r11 (Ohm/length_unit)"
$lineData[13,10] = "x11 (Ohm/length_unit)"
$lineData[13,11] = "r21 (Ohm/length_unit)"
$lineData[13,12] = "x21 (Ohm/length_unit)"
$lineData[13,13] = "r22 (Ohm/length_unit)"
$lineData[13,14] = "x22 (Ohm/length_unit)"
$lineData[13,15] = "r31 (Ohm/length_unit)"
$lineData[13,16] = "x31 (Ohm/length_unit)"
$lineData[13,17] = "r32 (Ohm/length_unit)"
$lineData[13,18] = "x32 (Ohm/length_unit)"
$lineData[13,19] = "r33 (Ohm/length_unit)"
$lineData[13,20] = "x33 (Ohm/length_unit)"
$lineData[13,21] = "b11 (uS/length_unit)"
$lineData[13,22] = "b21 (uS/length_unit)"
$lineData[13,23] = "b22 (uS/length_unit)"
$lineData[13,24] = "b31 (uS/length_unit)"
$lineData[13,25] = "b32 (uS/length_unit)"
$lineData[13,26] = "b33 (uS/length_unit)"
$lineData[14,0] = "End of Three-Phase Line with Full Data"
$lineData[16,0] = "Three-Phase Line with Sequential Data"
$lineData[17,0] = "ID"
$lineData[17,1] = "Status"
$lineData[17,2] = "Length"
$lineData[17,3] = "From1"
$lineData[17,4] = "From2"
$lineData[17,5] = "From3"
$lineData[17,6] = "To1"
$lineData[17,7] = "To2"
$lineData[17,8] = "To3"
$lineData[17,9] = "R0 (Ohm/length_unit)"
$lineData[17,10] = "X0 (Ohm/length_unit)"
$lineData[17,11] = "R1 (Ohm/length_unit)"
$lineData[17,12] = "X1 (Ohm/length_unit)"
$lineData[17,13] = "B0 (uS/length_unit)"
$lineData[17,14] = "B1 (uS/length_unit)"
$lineData[18,0] = "End of Three-Phase Line with Sequential Data"
$lineWs = $wb.Worksheets.Item("Line")
$lineWs.Range("A11:AA29").Value = $lineData
$styleSrc.Copy()
$lineWs.Range("A12:G12").PasteSpecial(-4122)
$styleSrc.Copy()
$lineWs.Range("A16:H16").PasteSpecial(-4122)
$styleSrc.Copy()
$lineWs.Range("A20:P20").PasteSpecial(-4122)
$styleSrc.Copy()
$lineWs.Range("A24:AA24").PasteSpecial(-4122)
$styleSrc.Copy()
$lineWs.Range("A28:O28").PasteSpecial(-4122)

# ===== Populate 'Transformer' sheet (dimension A11:AA25) =====
$transformerData = New-Object 'object[,]' 15,27
$transformerData[0,0] = "Positive-Sequence 2W Transformer"
$transformerData[1,0] = "ID"
$transformerData[1,1] = "Status"
$transformerData[1,2] = "From bus"
$transformerData[1,3] = "To bus"
$transformerData[1,4] = "R (pu)"
$transformerData[1,5] = "Xl (pu)"
$transformerData[1,6] = "Gmag (pu)"
$transformerData[1,7] = "Bmag (pu)"
$transformerData[1,8] = "Ratio W1 (pu)"
$transformerData[1,9] = "Ratio W2 (pu)"
$transformerData[1,10] = "Phase Shift (deg)"
$transformerData[2,0] = "End of Positive-Sequence 2W Transformer"
$transformerData[4,0] = "Positive-Sequence 3W Transformer"
$transformerData[5,0] = "ID"
$transformerData[5,1] = "Status"
$transformerData[5,2] = "Bus1"
$transformerData[5,3] = "Bus2"
$transformerData[5,4] = "Bus3"
$transformerData[5,5] = "R_12 (pu)"
$transformerData[5,6] = "Xl_12 (pu)"
$transformerData[5,7] = "R_23 (pu)"
$transformerData[5,8] = "Xl_23 (pu)"
$transformerData[5,9] = "R_31 (pu)"
$transformerData[5,10] = "Xl_31 (pu)"
$transformerData[5,11] = "Gmag (pu)"
$transformerData[5,12] = "Bmag (pu)"
$transformerData[5,13] = "Ratio W1 (pu)"
$transformerData[5,14] = "Ratio W2 (pu)"
$transformerData[5,15] = "Ratio W3 (pu)"
$transformerData[5,16] = "Phase Shift W1 (deg)"
$transformerData[5,17] = "Phase Shift W2 (deg)"
$transformerData[5,18] = "Phase Shift W3 (deg)"
$transformerData[6,0] = "End of Positive-Sequence 3W Transformer"
$transformerData[8,0] = "Multiphase 2W Transformer"
$transformerData[9,0] = "ID"
$transformerData[9,1] = "Status"
$transformerData[9,2] = "Number of phases"
$transformerData[9,3] = "Bus1_A"
$transformerData[9,4] = "Bus1_B"
$transformerData[9,5] = "Bus1_C"
$transformerData[9,6] = "V1 (kV)"
$transformerData[9,7] = "S_base1 (kVA)"
$transformerData[9,8] = "Conn. type1"
$transformerData[9,9] = "Bus2_A"
$transformerData[9,10] = "Bus2_B"
$transformerData[9,11] = "Bus2_C"
$transformerData[9,12] = "V2 (kV)"
$transformerData[9,13] = "S_base2 (kVA)"
$transformerData[9,14] = "Conn. type2"
$transformerData[9,15] = "Tap 1"
$transformerData[9,16] = "Tap 2"
$transformerData[9,17] = "Tap 3"
$transformerData[9,18] = "Lowest Tap"
$transformerData[9,19] = "Highest Tap"
$transformerData[9,20] = "Min Range (%)"
$transformerData[9,21] = "Max Range (%)"
$transformerData[9,22] = "X (pu)"
$transformerData[9,23] = "RW1 (pu)"
$transformerData[9,24] = "RW2"
$transformerData[10,0] = "End of Multiphase 2W Transformer"
$transformerData[12,0] = "Multiphase 2W Transformer with Mutual Impedance"
$transformerData[13,0] = "ID"
$transformerData[13,1] = "Status"
$transformerData[13,2] = "Number of phases"
$transformerData[13,3] = "Bus1_A"
$transformerData[13,4] = "Bus1_B"
$transformerData[13,5] = "Bus1_C"
$transformerData[13,6] = "V1 (kV)"
$transformerData[13,7] = "S_base1 (kVA)"
$transformerData[13,8] = "Conn. type1"
$transformerData[13,9] = "Bus2_A"
$transformerData[13,10] = "Bus2_B"
$transformerData[13,11] = "Bus2_C"
$transformerData[13,12] = "V2 (kV)"
$transformerData[13,13] = "S_base2 (kVA)"
$transformerData[13,14] = "Conn. type2"
$transformerData[13,15] = "Tap 1"
$transformerData[13,16] = "Tap 2"
$transformerData[13,17] = "Tap 3"
$transformerData[13,18] = "Lowest Tap"
$transformerData[13,19] = "Highest Tap"
$transformerData[13,20] = "Min Range (%)"
$transformerData[13,21] = "Max Range (%)"
$transformerData[13,22] = "Z0 leakage (pu)"
$transformerData[13,23] = "Z1 leakage (pu)"
$transformerData[13,24] = "X0/R0"
$transformerData[13,25] = "X1/R1"
$transformerData[13,26] = "No Load Loss (kW)"
$transformerData[14,0] = "End of Multiphase 2W Transformer with Mutual Impedance"
$transformerWs = $wb.Worksheets.Item("Transformer")
$transformerWs.Range("A11:AA25").Value = $transformerData
$styleSrc.Copy()
$transformerWs.Range("A12:K12").PasteSpecial(-4122)
$styleSrc.Copy()
$transformerWs.Range("A16:S16").PasteSpecial(-4122)
$styleSrc.Copy()
$transformerWs.Range("A20:Y20").PasteSpecial(-4122)
$styleSrc.Copy()
$transformerWs.Range("A24:AA24").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ===== Update the 'Bus' sheet: swap the A (bus id) and E (angle) values
# between specific row pairs so that phase A precedes phase C
# (except for the house_s1/house_s2 pair, which swaps the other way). =====
$busWs = $wb.Worksheets.Item("Bus")
$tmpA = $busWs.Range("A4").Value
$tmpE = $busWs.Range("E4").Value
$busWs.Range("A4").Value = $busWs.Range("A5").Value
$busWs.Range("E4").Value = $busWs.Range("E5").Value
$busWs.Range("A5").Value = $tmpA
$busWs.Range("E5").Value = $tmpE
$tmpA = $busWs.Range("A7").Value
$tmpE = $busWs.Range("E7").Value
$busWs.Range("A7").Value = $busWs.Range("A8").Value
$busWs.Range("E7").Value = $busWs.Range("E8").Value
$busWs.Range("A8").Value = $tmpA
$busWs.Range("E8").Value = $tmpE
$tmpA = $busWs.Range("A10").Value
$tmpE = $busWs.Range("E10").Value
$busWs.Range("A10").Value = $busWs.Range("A11").Value
$busWs.Range("E10").Value = $busWs.Range("E11").Value
$busWs.Range("A11").Value = $tmpA
$busWs.Range("E11").Value = $tmpE
$tmpA = $busWs.Range("A17").Value
$tmpE = $busWs.Range("E17").Value
$busWs.Range("A17").Value = $busWs.Range("A18").Value
$busWs.Range("E17").Value = $busWs.Range("E18").Value
$busWs.Range("A18").Value = $tmpA
$busWs.Range("E18").Value = $tmpE
$tmpA = $busWs.Range("A21").Value
$tmpE = $busWs.Range("E21").Value
$busWs.Range("A21").Value = $busWs.Range("A22").Value
$busWs.Range("E21").Value = $busWs.Range("E22").Value
$busWs.Range("A22").Value = $tmpA
$busWs.Range("E22").Value = $tmpE
$tmpA = $busWs.Range("A24").Value
$tmpE = $busWs.Range("E24").Value
$busWs.Range("A24").Value = $busWs.Range("A25").Value
$busWs.Range("E24").Value = $busWs.Range("E25").Value
$busWs.Range("A25").Value = $tmpA
$busWs.Range("E25").Value = $tmpE
$tmpA = $busWs.Range("A27").Value
$tmpE = $busWs.Range("E27").Value
$busWs.Range("A27").Value = $busWs.Range("A28").Value
$busWs.Range("E27").Value = $busWs.Range("E28").Value
$busWs.Range("A28").Value = $tmpA
$busWs.Range("E28").Value = $tmpE
$tmpA = $busWs.Range("A30").Value
$tmpE = $busWs.Range("E30").Value
$busWs.Range("A30").Value = $busWs.Range("A31").Value
$busWs.Range("E30").Value = $busWs.Range("E31").Value
$busWs.Range("A31").Value = $tmpA
$busWs.Range("E31").Value = $tmpE
$tmpA = $busWs.Range("A32").Value
$tmpE = $busWs.Range("E32").Value
$busWs.Range("A32").Value = $busWs.Range("A33").Value
$busWs.Range("E32").Value = $busWs.Range("E33").Value
$busWs.Range("A33").Value = $tmpA
$busWs.Range("E33").Value = $tmpE
$tmpA = $busWs.Range("A35").Value
$tmpE = $busWs.Range("E35").Value
$busWs.Range("A35").Value = $busWs.Range("A36").Value
$busWs.Range("E35").Value = $busWs.Range("E36").Value
$busWs.Range("A36").Value = $tmpA
$busWs.Range("E36").Value = $tmpE
$tmpA = $busWs.Range("A38").Value
$tmpE = $busWs.Range("E38").Value
$busWs.Range("A38").Value = $busWs.Range("A39").Value
$busWs.Range("E38").Value = $busWs.Range("E39").Value
$busWs.Range("A39").Value = $tmpA
$busWs.Range("E39").Value = $tmpE
$tmpA = $busWs.Range("A40").Value
$tmpE = $busWs.Range("E40").Value
$busWs.Range("A40").Value = $busWs.Range("A41").Value
$busWs.Range("E40").Value = $busWs.Range("E41").Value
$busWs.Range("A41").Value = $tmpA
$busWs.Range("E41").Value = $tmpE
$tmpA = $busWs.Range("A43").Value
$tmpE = $busWs.Range("E43").Value
$busWs.Range("A43").Value = $busWs.Range("A44").Value
$busWs.Range("E43").Value = $busWs.Range("E44").Value
$busWs.Range("A44").Value = $tmpA
$busWs.Range("E44").Value = $tmpE
$tmpA = $busWs.Range("A46").Value
$tmpE = $busWs.Range("E46").Value
$busWs.Range("A46").Value = $busWs.Range("A47").Value
$busWs.Range("E46").Value = $busWs.Range("E47").Value
$busWs.Range("A47").Value = $tmpA
$busWs.Range("E47").Value = $tmpE
$tmpA = $busWs.Range("A49").Value
$tmpE = $busWs.Range("E49").Value
$busWs.Range("A49").Value = $busWs.Range("A50").Value
$busWs.Range("E49").Value = $busWs.Range("E50").Value
$busWs.Range("A50").Value = $tmpA
$busWs.Range("E50").Value = $tmpE
$tmpA = $busWs.Range("A53").Value
$tmpE = $busWs.Range("E53").Value
$busWs.Range("A53").Value = $busWs.Range("A54").Value
$busWs.Range("E53").Value = $busWs.Range("E54").Value
$busWs.Range("A54").Value = $tmpA
$busWs.Range("E54").Value = $tmpE
